# The worksheet ("Sheet1") ships with sheetProtection enabled (password-hashed).
# All data cells are locked by default, so a direct Range.Value write throws
# "protected sheet". We avoid Unprotect()/Protect() (which would strip/alter
# the original <sheetProtection> element and its legacy password hash) by
# temporarily unlocking just column E, writing the new "Percent Change"
# values, and then restoring the original cell formatting/lock state by
# copying the format from the untouched column D (which carries the same
# style) back onto column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Percent Change" (column E) values for rows 2-37, taken from the diff.
$values = @{
    2  = 0.008813617903741466
    3  = 0.01486011057128511
    4  = -0.0416005080978088
    5  = 0.005669606664236326
    6  = 0.008021390374331583
    7  = -0.008640838586302158
    8  = -0.001764826595526348
    9  = 0.01373795366003683
    10 = -0.01048667011299009
    11 = -0.01701323251417775
    12 = -0.01104926476696844
    13 = -0.0161094751429065
    14 = 0.003232570330246487
    15 = 0.00458520845716226
    16 = 0.001880091915604787
    17 = -0.007232767232767223
    18 = -0.0139640298289222
    19 = 0.006485194982507103
    20 = -0.006491297692612807
    21 = 0.002361852802243991
    22 = -0.005921159357594741
    23 = 0.0132180667794275
    24 = 0.02404612930928729
    25 = -0.03025189827648556
    26 = -0.03413173652694623
    27 = 0.007493875198155475
    28 = 0.002640752086194187
    29 = 0.008175466397752684
    30 = 0.0003122853038535212
    31 = -0.01578616807981259
    32 = 0.004425668879501021
    33 = -0.014993752603082
    34 = 0.003412470140886414
    35 = 0.00117332638028822
    36 = 0.008555713994703673
    37 = -0.001484514705073936
}

# Temporarily unlock column E (rows 2-37) so values can be written while the
# sheet stays protected.
$editRange = $ws.Range("E2:E37")
$editRange.Locked = $false

foreach ($row in $values.Keys) {
    $ws.Range("E$row").Value = $values[$row]
}

# Restore the original cell formatting/protection (style "s=1") by copying
# the format from column D, which was never touched and still carries it.
$ws.Range("D2:D37").Copy()
$editRange.PasteSpecial(-4122)  # xlPasteFormats
